$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet from "Data" to "Summary" ---
$ws.Name = "Summary"

# --- Add the new bold+underline "title_" cell style (used by the new
#     "Source Type: Statistical Institution" heading) ---
$titleUnderline = $wb.Styles.Add("title_")
$titleUnderline.Font.Bold = $true
$titleUnderline.Font.Underline = $true

# --- Move the old "Micro / SMEs / MSMEs" header row (was row 5) down to
#     row 9, and the "Enterprises (% of total)" row (was row 6) down to
#     row 10, clearing out the old, now-vacated cells ---
$oldMicro = $ws.Range("B5").Value()
$oldSMEs = $ws.Range("C5").Value()
$oldMSMEs = $ws.Range("D5").Value()
$oldEnterprises = $ws.Range("A6").Value()

$ws.Range("B5:D5").Clear()
$ws.Range("A6").Clear()

$ws.Range("B9").Value = $oldMicro
$ws.Range("B9").Style = "title"
$ws.Range("C9").Value = $oldSMEs
$ws.Range("C9").Style = "title"
$ws.Range("D9").Value = $oldMSMEs
$ws.Range("D9").Style = "title"

$ws.Range("A10").Value = $oldEnterprises
$ws.Range("A10").Style = "title"

# --- New "Source Type: Statistical Institution" heading at A7 ---
$ws.Range("A7").Value = "Source Type: Statistical Institution"
$ws.Range("A7").Style = "title_"

# --- New MSME participation value at D10 (stored as text, Normal style) ---
$ws.Range("D10").Value = "'99.5"
$ws.Range("D10").Style = "Normal"

# --- New source citation rows ---
$ws.Range("A11").Value = "Source: BSC, 2009"
$ws.Range("A11").Style = "source"

$ws.Range("A19").Value = "BSC"
$ws.Range("A19").Style = "title"

$ws.Range("A20").Value = "Bureau of Statistics and Census (BSC) Libya, ""المنشأت التحويلية الصغيرة 2009"" and ""المنشأت الصناعية التحويلية الكبيرة 2009"". Available at http://bsc.ly/?P=5&sec_Id=12&dep_Id=4"
$ws.Range("A20").Style = "source"

Write-Output "applied edits"
